$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4180193265783032
$ws.Range("D2").Value = 0.06747806412783142
$ws.Range("E2").Value = 0.1765666249016249
$ws.Range("F2").Value = 1.54402593499843
$ws.Range("G2").Value = 1.454626884411056
$ws.Range("H2").Value = 1.314481704399967
$ws.Range("K2").Value = 2.225591597444293
$ws.Range("L2").Value = 0.1540696821842005
$ws.Range("C3").Value = 0.4095275913296064
$ws.Range("D3").Value = 0.06844969275959478
$ws.Range("E3").Value = 0.1725842348115449
$ws.Range("F3").Value = 1.540913468417642
$ws.Range("G3").Value = 1.453202398683601
$ws.Range("H3").Value = 1.32210160019666
$ws.Range("K3").Value = 2.028806800577968
$ws.Range("L3").Value = 0.150242226801268
$ws.Range("C4").Value = 0.4045550471834929
$ws.Range("D4").Value = 0.06907795365981784
$ws.Range("E4").Value = 0.1702440526323663
$ws.Range("F4").Value = 1.540311051451809
$ws.Range("G4").Value = 1.453684529586752
$ws.Range("H4").Value = 1.32770002083565
$ws.Range("K4").Value = 1.908458429628411
$ws.Range("L4").Value = 0.1479862963206315
$ws.Range("C5").Value = 0.4025892243720648
$ws.Range("D5").Value = 0.06934193376535625
$ws.Range("E5").Value = 0.1693167216340576
$ws.Range("F5").Value = 1.540392972027448
$ws.Range("G5").Value = 1.45422015581407
$ws.Range("H5").Value = 1.33021193112647
$ws.Range("K5").Value = 1.85953615494185
$ws.Range("L5").Value = 0.1470905554261748
$ws.Range("C6").Value = 0.4022664517434293
$ws.Range("D6").Value = 0.0693862479642231
$ws.Range("E6").Value = 0.1691643254512059
$ws.Range("F6").Value = 1.54042630119767
$ws.Range("G6").Value = 1.454329521300849
$ws.Range("H6").Value = 1.330642931416179
$ws.Range("K6").Value = 1.851419949629076
$ws.Range("L6").Value = 0.1469432388523515
$ws.Range("C7").Value = 0.4045282905583463
$ws.Range("D7").Value = 0.06908148157330807
$ws.Range("E7").Value = 0.1702314399097986
$ws.Range("F7").Value = 1.540310832802177
$ws.Range("G7").Value = 1.453690382691377
$ws.Range("H7").Value = 1.32773296506079
$ws.Range("K7").Value = 1.907798157581226
$ws.Range("L7").Value = 0.1479741207496232
$ws.Range("C8").Value = 0.4150411593057868
$ws.Range("D8").Value = 0.06780650291752721
$ws.Range("E8").Value = 0.175171639472957
$ws.Range("F8").Value = 1.542680139536415
$ws.Range("G8").Value = 1.453852903271937
$ws.Range("H8").Value = 1.316917675562138
$ws.Range("K8").Value = 2.157641080937026
$ws.Range("L8").Value = 0.1527303717327655
$ws.Range("C9").Value = 0.4375821139846607
$ws.Range("D9").Value = 0.06555792139741179
$ws.Range("E9").Value = 0.1856977882698132
$ws.Range("F9").Value = 1.557787587333507
$ws.Range("G9").Value = 1.465030421495015
$ws.Range("H9").Value = 1.303042935183726
$ws.Range("K9").Value = 2.651389842978347
$ws.Range("L9").Value = 0.1628097635934438
$ws.Range("C10").Value = 0.4553326513109539
$ws.Range("D10").Value = 0.06405972410548344
$ws.Range("E10").Value = 0.1939507371160758
$ws.Range("F10").Value = 1.575376121500327
$ws.Range("G10").Value = 1.47999532905439
$ws.Range("H10").Value = 1.29737151372251
$ws.Range("K10").Value = 3.016529332124094
$ws.Range("L10").Value = 0.1706823420872041
$ws.Range("C11").Value = 0.4636698452827943
$ws.Range("D11").Value = 0.06341166844713086
$ws.Range("E11").Value = 0.197819881659079
$ws.Range("F11").Value = 1.584811108444555
$ws.Range("G11").Value = 1.488298588482905
$ws.Range("H11").Value = 1.295784689372937
$ws.Range("K11").Value = 3.183175130614529
$ws.Range("L11").Value = 0.1743671849165338
$ws.Range("C12").Value = 0.4668649128249456
$ws.Range("D12").Value = 0.06317109372472274
$ws.Range("E12").Value = 0.1993016740729558
$ws.Range("F12").Value = 1.588592015970292
$ws.Range("G12").Value = 1.491660205753647
$ws.Range("H12").Value = 1.295327522451771
$ws.Range("K12").Value = 3.246358148612956
$ws.Range("L12").Value = 0.1757775727584772
$ws.Range("C13").Value = 0.46617510628829
$ws.Range("D13").Value = 0.06322269073042008
$ws.Range("E13").Value = 0.1989818019746394
$ws.Range("F13").Value = 1.587768445042713
$ws.Range("G13").Value = 1.490926517996542
$ws.Range("H13").Value = 1.295419574434135
$ws.Range("K13").Value = 3.232747092081922
$ws.Range("L13").Value = 0.1754731505629934
$ws.Range("C14").Value = 0.4639319434819811
$ws.Range("D14").Value = 0.06339177933098128
$ws.Range("E14").Value = 0.197941455610767
$ws.Range("F14").Value = 1.58511798479897
$ws.Range("G14").Value = 1.488570781558622
$ws.Range("H14").Value = 1.295744192146685
$ws.Range("K14").Value = 3.188371674982989
$ws.Range("L14").Value = 0.1744829166687225
$ws.Range("C15").Value = 0.4625628908069928
$ws.Range("D15").Value = 0.06349598044508298
$ws.Range("E15").Value = 0.1973063828616546
$ws.Range("F15").Value = 1.583521656318638
$ws.Range("G15").Value = 1.487156198722801
$ws.Range("H15").Value = 1.295961775004542
$ws.Range("K15").Value = 3.161200591146724
$ws.Range("L15").Value = 0.1738783293330215
$ws.Range("C16").Value = 0.4547930948323824
$ws.Range("D16").Value = 0.06410275108528296
$ws.Range("E16").Value = 0.1937002002543409
$ws.Range("F16").Value = 1.574788537657057
$ws.Range("G16").Value = 1.479482982392341
$ws.Range("H16").Value = 1.297495279062474
$ws.Range("K16").Value = 3.005649536047031
$ws.Range("L16").Value = 0.1704436236088043
$ws.Range("C17").Value = 0.4500939363948646
$ws.Range("D17").Value = 0.06448357123230686
$ws.Range("E17").Value = 0.191517419586809
$ws.Range("F17").Value = 1.569799644523656
$ws.Range("G17").Value = 1.475160401881283
$ws.Range("H17").Value = 1.298691049733549
$ws.Range("K17").Value = 2.910362747974489
$ws.Range("L17").Value = 0.1683631571192024
$ws.Range("C18").Value = 0.4474157880852658
$ws.Range("D18").Value = 0.0647057593976097
$ws.Range("E18").Value = 0.190272743644627
$ws.Range("F18").Value = 1.567065014197624
$ws.Range("G18").Value = 1.472814808108126
$ws.Range("H18").Value = 1.299472248509602
$ws.Range("K18").Value = 2.855607323433333
$ws.Range("L18").Value = 0.167176268450433
$ws.Range("C19").Value = 0.4465132457832794
$ws.Range("D19").Value = 0.0647815291447511
$ws.Range("E19").Value = 0.1898531694184058
$ws.Range("F19").Value = 1.566162218335762
$ws.Range("G19").Value = 1.472044716193608
$ws.Range("H19").Value = 1.299752766663971
$ws.Range("K19").Value = 2.83707684304045
$ws.Range("L19").Value = 0.1667760775362126
$ws.Range("C20").Value = 0.4505916136081964
$ws.Range("D20").Value = 0.06444270610668035
$ws.Range("E20").Value = 0.191748661527285
$ws.Range("F20").Value = 1.570316750684384
$ws.Range("G20").Value = 1.475605975761113
$ws.Range("H20").Value = 1.298554082349938
$ws.Range("K20").Value = 2.920500910074736
$ws.Range("L20").Value = 0.1685836173138995
$ws.Range("C21").Value = 0.4645897829363435
$ws.Range("D21").Value = 0.06334198270161551
$ws.Range("E21").Value = 0.1982465782331815
$ws.Range("F21").Value = 1.585890826753385
$ws.Range("G21").Value = 1.489256801296079
$ws.Range("H21").Value = 1.295644935935172
$ws.Range("K21").Value = 3.201403697406647
$ws.Range("L21").Value = 0.1747733637791669
$ws.Range("C22").Value = 0.4739597228576145
$ws.Range("D22").Value = 0.06265075643934281
$ws.Range("E22").Value = 0.2025903392651642
$ws.Range("F22").Value = 1.597283140258753
$ws.Range("G22").Value = 1.499446375236147
$ws.Range("H22").Value = 1.294581835243605
$ws.Range("K22").Value = 3.385443825962966
$ws.Range("L22").Value = 0.1789063161437667
$ws.Range("C23").Value = 0.4689384832013559
$ws.Range("D23").Value = 0.06301709494217178
$ws.Range("E23").Value = 0.2002630773779259
$ws.Range("F23").Value = 1.591091160133828
$ws.Range("G23").Value = 1.49389123068903
$ws.Range("H23").Value = 1.295072231794649
$ws.Range("K23").Value = 3.287176613206213
$ws.Range("L23").Value = 0.1766924247657613
$ws.Range("C24").Value = 0.4503665406828077
$ws.Range("D24").Value = 0.06446117110906258
$ws.Range("E24").Value = 0.1916440851790568
$ws.Range("F24").Value = 1.570082551151202
$ws.Range("G24").Value = 1.475404097423251
$ws.Range("H24").Value = 1.298615713411834
$ws.Range("K24").Value = 2.915917365314613
$ws.Range("L24").Value = 0.1684839186155358
$ws.Range("C25").Value = 0.4312764409112901
$ws.Range("D25").Value = 0.06613924157371365
$ws.Range("E25").Value = 0.1827595815844916
$ws.Range("F25").Value = 1.552569179461997
$ws.Range("G25").Value = 1.460830401108353
$ws.Range("H25").Value = 1.306005983682951
$ws.Range("K25").Value = 2.517404092411255
$ws.Range("L25").Value = 0.1600015945912929
